$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'34.191.08"
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.94%  '

$c = $ws.Range("D3")
$c.Value = "'1.815.13"
$c.Style = "Normal"
$ws.Range("E3").Value = '  +1.41%  '

$ws.Range("E4").Value = '  +0.01%  '

$c = $ws.Range("D5")
$c.Value = "'225.04"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.18%  '

$c = $ws.Range("D6")
$c.Value = "'0.556"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.50%  '

$ws.Range("E7").Value = '  -0.02%  '

$c = $ws.Range("D8")
$c.Value = "'32.01"
$c.Style = "Normal"
$ws.Range("E8").Value = '  -3.74%  '

$ws.Range("E9").Value = '  +3.41%  '

$c = $ws.Range("D10")
$c.Value = "'0.0733"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +10.30%  '

$ws.Range("E11").Value = '  -0.07%  '

$c = $ws.Range("D12")
$c.Value = "'2.078.69"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.53%  '

$c = $ws.Range("D13")
$c.Value = "'11.05"
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.47%  '

$c = $ws.Range("D14")
$c.Value = "'1.808.74"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.99%  '

$c = $ws.Range("D15")
$c.Value = "'0.642"
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.70%  '

$c = $ws.Range("D16")
$c.Value = "'34.206.58"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.72%  '

$c = $ws.Range("D17")
$c.Value = "'4.33"
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.55%  '

$c = $ws.Range("D18")
$c.Value = "'69.60"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.40%  '

$c = $ws.Range("D19")
$c.Value = "'249.74"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -2.80%  '

$c = $ws.Range("D20")
$c.Value = "'0.0₃0801"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +6.63%  '

$c = $ws.Range("D21")
$c.Value = "'11.04"
$c.Style = "Normal"
$ws.Range("E21").Value = '  +5.44%  '

$ws.Range("E22").Value = '  -0.16%  '

$c = $ws.Range("D23")
$c.Value = "'4.24"
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.60%  '

$ws.Range("E24").Value = '  +0.74%  '

$c = $ws.Range("D25")
$c.Value = "'160.50"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.78%  '

$c = $ws.Range("D26")
$c.Value = "'16.67"
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.90%  '

$c = $ws.Range("D27")
$c.Value = "'7.22"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.06%  '

$ws.Range("E28").Value = '  +0.43%  '

$ws.Range("E29").Value = '  +0.07%  '

$c = $ws.Range("D30")
$c.Value = "'0.0533"
$c.Style = "Normal"
$ws.Range("E30").Value = '  +3.41%  '

$c = $ws.Range("D31")
$c.Value = "'3.79"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.29%  '

$ws.Range("E32").Value = '  +2.10%  '

$c = $ws.Range("D33")
$c.Value = "'3.60"
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.40%  '

$ws.Range("E34").Value = '  -1.12%  '

$c = $ws.Range("D35")
$c.Value = "'1.430.63"
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.74%  '

$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D36")
$c.Value = "'1.06"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.61%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D37")
$c.Value = "'0.641"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.54%  '

$c = $ws.Range("D38")
$c.Value = "'0.0190"
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.50%  '

$c = $ws.Range("D39")
$c.Value = "'0.961"
$c.Style = "Normal"
$ws.Range("E39").Value = '  +7.30%  '

$c = $ws.Range("D40")
$c.Value = "'81.69"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -3.02%  '

$c = $ws.Range("D41")
$c.Value = "'2.75"
$c.Style = "Normal"
$ws.Range("E41").Value = '  -3.96%  '

$c = $ws.Range("D42")
$c.Value = "'2.35"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.13%  '

$c = $ws.Range("D43")
$c.Value = "'2.16"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +3.59%  '

$c = $ws.Range("D44")
$c.Value = "'6.01"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.36%  '

$ws.Range("E45").Value = '  -1.55%  '

$c = $ws.Range("D46")
$c.Value = "'1.972.92"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.18%  '

$c = $ws.Range("D47")
$c.Value = "'106.91"
$c.Style = "Normal"
$ws.Range("E47").Value = '  +7.65%  '

$ws.Range("E48").Value = '  -1.23%  '

$ws.Range("E49").Value = '  -0.17%  '

$c = $ws.Range("D50")
$c.Value = "'11.89"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -6.52%  '

$ws.Range("E51").Value = '  +6.86%  '
